$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42626
$ws.Range("A5").NumberFormat = "d-mmm"

$ws.Range("B5").Value = 0.89583333333333337
$ws.Range("B5").NumberFormat = "h:mm"

$ws.Range("C5").Value = 0.95833333333333337
$ws.Range("C5").NumberFormat = "h:mm"

$ws.Range("D5").Value = "Estudo e implantação do sessionStorage + login"

$ws.Range("D6").Select()
